$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header stays the same label, new dataset follows below it.
$ws.Range("A1").Value = "HK_G_acc_G"

$values = @(
    70.652173913043484,
    70.380434782608688,
    70.380434782608688,
    70.108695652173907,
    69.83695652173914,
    70.108695652173907,
    71.739130434782609,
    72.282608695652172,
    72.554347826086953,
    72.282608695652172,
    70.652173913043484,
    71.739130434782609,
    72.554347826086953,
    72.554347826086953,
    72.554347826086953,
    71.467391304347828,
    71.739130434782609,
    71.739130434782609,
    71.739130434782609,
    72.282608695652172,
    72.554347826086953,
    69.021739130434781,
    70.380434782608688,
    69.565217391304344,
    71.195652173913047,
    71.739130434782609,
    72.282608695652172,
    73.369565217391312,
    73.097826086956516,
    73.097826086956516,
    69.565217391304344,
    69.83695652173914,
    70.923913043478265,
    70.652173913043484,
    70.108695652173907,
    71.195652173913047,
    69.021739130434781,
    69.021739130434781,
    69.565217391304344,
    72.282608695652172,
    72.282608695652172,
    72.282608695652172,
    71.739130434782609,
    69.565217391304344,
    70.380434782608688,
    69.565217391304344,
    70.380434782608688,
    69.83695652173914
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
